$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.716.14"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.22%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.872.38"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.04%  "
$ws.Range("E4").Value = "  +0.45%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.006"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.38%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4586"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.64%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3845"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.45%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07853"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.11%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9914"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.14%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.56"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.98%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.931.94"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.86%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.958"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.85%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.674"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.06%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06960"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.65%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.25"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.09%  "
$ws.Range("E17").Value = "  +0.34%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001003"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.97%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.75"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.45%  "
$ws.Range("E20").Value = "  +0.37%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "28.715.66"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.263"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.57%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.119"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.58%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.119.05"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.79%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.13"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.21"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.742"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "118.71"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.77%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.933"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.13%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09296"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.60%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9111"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.92%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.286"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.33%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.330"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.80%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.306"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.81%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05738"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.55%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.144"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.58%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02071"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.19%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.666"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.73%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5614"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.22%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1781"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.40%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.822"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.18%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.07170"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.09%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "11.65"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.45%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5261"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.124"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.111"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.89%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.821"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.81%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "113.08"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.52%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.415"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.33%  "
$ws.Range("E51").Value = "  +0.36%  "
